# Auto-generated edit script: updates the cryptos price table
# to the latest snapshot values (coin list shifted by one rank,
# price/volume figures refreshed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.017.08"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3
$ws.Range("D3").Value = "2.116.31"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5192"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4452"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09382"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.183"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.584"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.13%  "

# Row 14
$ws.Range("D14").Value = "2.114.71"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.914"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.85%  "

# Row 17
$ws.Range("E17").Value = "  -0.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06694"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.301"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.94%  "

# Row 22
$ws.Range("E22").Value = "  +0.18%  "

# Row 23
$ws.Range("D23").Value = "30.034.69"
$ws.Range("E23").Value = "  -0.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.318"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.536"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "134.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.153"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.67%  "

# Row 31
$ws.Range("B31").Value = "ARBITRUM"
$ws.Range("C31").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.790"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.73%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1056"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.63%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.253"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.644"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.53%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.968"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.33%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.95%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02615"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.74%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "

# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7083"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.39%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.334"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2237"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.51%  "

# Row 43
$ws.Range("B43").Value = "Decentraland"
$ws.Range("C43").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6840"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.11%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.36%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.375"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.32%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.269"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.19%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000358"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.627"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.38%  "

# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.225"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "

